# Updates the practice sheet from 2024-01-15 (Monday) to 2024-01-16 (Tuesday)
# and refreshes every "three-digit x one-digit" multiplication problem in the
# answer table with the new day's values.

$d = $word.ActiveDocument

# --- Title line: date + weekday -------------------------------------------
$d.Content.Find.Execute("2024-01-15 Monday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2024-01-16 Tuesday", 2) | Out-Null

# --- Row 1 of the table: cell contents are reshuffled -----------------------
# old: 983x2=1966 | 257x5=1285 | 584x7=4088 | 742x6=4452 | 280x5=1400
# new: 152x5=760  | 685x6=4110 | 141x4=564  | 120x6=720  | 742x6=4452
$t = $d.Tables.Item(1)
$row1 = $t.Rows.Item(1)
$row1.Cells.Item(1).Range.Text = "152×5=760"
$row1.Cells.Item(2).Range.Text = "685×6=4110"
$row1.Cells.Item(3).Range.Text = "141×4=564"
$row1.Cells.Item(4).Range.Text = "120×6=720"
$row1.Cells.Item(5).Range.Text = "742×6=4452"

# --- Remaining answer rows: straightforward value swaps ---------------------
$replacements = @(
    @("240×3=720", "790×5=3950"),
    @("786×6=4716", "980×7=6860"),
    @("896×2=1792", "153×5=765"),
    @("400×3=1200", "572×4=2288"),
    @("820×7=5740", "709×4=2836"),
    @("962×3=2886", "441×3=1323"),
    @("762×5=3810", "905×6=5430"),
    @("856×9=7704", "317×9=2853"),
    @("636×7=4452", "560×4=2240"),
    @("188×9=1692", "247×6=1482"),
    @("211×3=633", "910×8=7280"),
    @("420×9=3780", "481×6=2886"),
    @("364×2=728", "779×4=3116"),
    @("667×4=2668", "968×8=7744"),
    @("383×7=2681", "706×4=2824"),
    @("937×6=5622", "808×5=4040"),
    @("125×4=500", "175×8=1400"),
    @("421×7=2947", "256×6=1536"),
    @("863×4=3452", "882×9=7938"),
    @("313×3=939", "138×9=1242")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "done"
